$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loads")

# Re-lay the "loads" table out with two new leading data columns
# (v_nom_kv, s_base_mva) and two new trailing columns (g_shunt_pu,
# b_shunt_pu). Existing columns (v_nom_pu, p_nom_mw, q_nom_mvar, bus_idx)
# slide from B:E to D:G. Cells are written directly (not via a column
# insert) so per-column formatting stays anchored to its original column.

# Row 1 - headers
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "v_nom_kv"
$ws.Range("C1").Value = "s_base_mva"
$ws.Range("D1").Value = "v_nom_pu"
$ws.Range("E1").Value = "p_nom_mw"
$ws.Range("F1").Value = "q_nom_mvar"
$ws.Range("G1").Value = "bus_idx"
$ws.Range("H1").Value = "g_shunt_pu"
$ws.Range("I1").Value = "b_shunt_pu"

# Row 2 - Load 1
$ws.Range("B2").Value = 22
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 50
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0

# Row 3 - Load 2
$ws.Range("B3").Value = 132
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0

# Row 4 - Load 3
$ws.Range("B4").Value = 132
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0

# Selection on this sheet, and make it the active tab (was "trafos").
[void]$ws.Range("I5").Select()
$ws.Activate()
